{"js": "// 1. Update the letter date from September 19, 2025 to September 21, 2025.\nconst dateResults = context.document.body.search(\"September 19, 2025\", { matchCase: true });\ndateResults.load(\"items\");\nawait context.sync();\nif (dateResults.items.length > 0) {\n  dateResults.items[0].insertText(\"September 21, 2025\", Word.InsertLocation.replace);\n}\n\n// 2. Split the mailing-address line \"231 High Street, Palo Alto CA 94309\"\n// (the homeowner address block, not the PROPERTY ADDRESS table cell) into\n// two paragraphs: \"231 High Street\" and a new \"Palo Alto, CA 94309\" line.\nconst addressResults = context.document.body.search(\"231 High Street, Palo Alto CA 94309\", { matchCase: true });\naddressResults.load(\"items\");\nawait context.sync();\n\nfor (let i = 0; i < addressResults.items.length; i++) {\n  const candidate = addressResults.items[i];\n  const table = candidate.parentTableOrNullObject;\n  table.load(\"isNullObject\");\n  await context.sync();\n  if (table.isNullObject) {\n    const addressPara = candidate.paragraphs.getFirst();\n    addressPara.insertParagraph(\"Palo Alto, CA 94309\", Word.InsertLocation.after);\n    candidate.insertText(\"231 High Street\", Word.InsertLocation.replace);\n    break;\n  }\n}\nawait context.sync();\n\n// 3. Remove the blank \"No Spacing\" paragraph that sits right after the\n// \"Board of Directors\" line.\nconst boardResults = context.document.body.search(\"Board of Directors\", { matchCase: true });\nboardResults.load(\"items\");\nawait context.sync();\n\nif (boardResults.items.length > 0) {\n  const boardPara = boardResults.items[0].paragraphs.getFirst();\n  const nextPara = boardPara.getNext();\n  nextPara.load(\"text\");\n  await context.sync();\n  if (nextPara.text === \"\") {\n    nextPara.delete();\n  }\n}\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# 1. Update the letter date from September 19, 2025 to September 21, 2025.\n$find = $d.Content.Find\n$find.Text = \"September 19, 2025\"\n$find.Replacement.Text = \"September 21, 2025\"\n$find.Execute(\"September 19, 2025\", $false, $false, $false, $false, $false, $true, 1, $false, \"September 21, 2025\", 2)\n\n# 2. Split the mailing-address line \"231 High Street, Palo Alto CA 94309\"\n# (the homeowner address block, not the PROPERTY ADDRESS table cell) into\n# two paragraphs: \"231 High Street\" and a new \"Palo Alto, CA 94309\" line.\n$rng = $d.Content\n$find2 = $rng.Find\n$find2.Text = \"231 High Street, Palo Alto CA 94309\"\n$find2.Forward = $true\n$find2.Wrap = 0\n$guard = 0\nwhile ($find2.Execute() -and $guard -lt 20) {\n    $guard = $guard + 1\n    if ($rng.Information(12) -eq $false) {\n        $para = $rng.Paragraphs(1)\n        $rng.Text = \"231 High Street\"\n        $rng.InsertParagraphAfter()\n        $newPara = $para.Next()\n        $newPara.Range.Text = \"Palo Alto, CA 94309\"\n        break\n    }\n    $rng.Collapse(0)\n}\n\n# 3. Remove the blank \"No Spacing\" paragraph that sits right after the\n# \"Board of Directors\" line.\n$rng2 = $d.Content\n$find3 = $rng2.Find\n$find3.Text = \"Board of Directors\"\n$found3 = $find3.Execute()\nif ($found3) {\n    $boardPara = $rng2.Paragraphs(1)\n    $nextPara = $boardPara.Next()\n    if ($nextPara.Range.Text.Trim() -eq \"\") {\n        $nextPara.Range.Delete()\n    }\n}\n"}
